$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value (as text, matching the source data which
# stores every Coin/Link/Price/Volume cell as a string literal).
$edits = @(
    @{Cell='D2'; Value='34.855.29'}
    @{Cell='D3'; Value='1.839.14'}
    @{Cell='E3'; Value='  +1.57%  '}
    @{Cell='E4'; Value='  -0.02%  '}
    @{Cell='D5'; Value='231.79'}
    @{Cell='E5'; Value='  -0.47%  '}
    @{Cell='D6'; Value='0.620'}
    @{Cell='E6'; Value='  +1.05%  '}
    @{Cell='E7'; Value='  -0.06%  '}
    @{Cell='D8'; Value='39.98'}
    @{Cell='E8'; Value='  -0.90%  '}
    @{Cell='D9'; Value='0.329'}
    @{Cell='E9'; Value='  +1.29%  '}
    @{Cell='D10'; Value='0.0687'}
    @{Cell='E10'; Value='  +0.42%  '}
    @{Cell='D11'; Value='0.0982'}
    @{Cell='E11'; Value='  -1.52%  '}
    @{Cell='D12'; Value='2.105.64'}
    @{Cell='E12'; Value='  +1.50%  '}
    @{Cell='D13'; Value='11.42'}
    @{Cell='E13'; Value='  +3.61%  '}
    @{Cell='D14'; Value='1.855.11'}
    @{Cell='E14'; Value='  +2.52%  '}
    @{Cell='E15'; Value='  +1.42%  '}
    @{Cell='D16'; Value='4.65'}
    @{Cell='E16'; Value='  +0.14%  '}
    @{Cell='D17'; Value='34.903.34'}
    @{Cell='E17'; Value='  -0.48%  '}
    @{Cell='D18'; Value='69.84'}
    @{Cell='E18'; Value='  +0.31%  '}
    @{Cell='D19'; Value='0.0₃0788'}
    @{Cell='E19'; Value='  -0.23%  '}
    @{Cell='D20'; Value='240.41'}
    @{Cell='E20'; Value='  +0.87%  '}
    @{Cell='D21'; Value='12.21'}
    @{Cell='E21'; Value='  +2.33%  '}
    @{Cell='E22'; Value='  -0.31%  '}
    @{Cell='E23'; Value='  +0.01%  '}
    @{Cell='E24'; Value='  +0.90%  '}
    @{Cell='D25'; Value='171.24'}
    @{Cell='E25'; Value='  -0.59%  '}
    @{Cell='D26'; Value='7.80'}
    @{Cell='E26'; Value='  -0.49%  '}
    @{Cell='D27'; Value='17.46'}
    @{Cell='E27'; Value='  -0.29%  '}
    @{Cell='E28'; Value='  +2.28%  '}
    @{Cell='D29'; Value='1.51'}
    @{Cell='E29'; Value='  -5.24%  '}
    @{Cell='E30'; Value='  -0.05%  '}
    @{Cell='D31'; Value='0.0552'}
    @{Cell='E31'; Value='  -0.27%  '}
    @{Cell='E32'; Value='  -4.60%  '}
    @{Cell='D33'; Value='3.96'}
    @{Cell='E33'; Value='  -1.58%  '}
    @{Cell='D34'; Value='1.90'}
    @{Cell='E34'; Value='  +7.22%  '}
    @{Cell='D35'; Value='1.23'}
    @{Cell='E35'; Value='  +7.47%  '}
    @{Cell='D36'; Value='1.45'}
    @{Cell='E36'; Value='  +13.11%  '}
    @{Cell='D37'; Value='0.696'}
    @{Cell='E37'; Value='  +2.37%  '}
    @{Cell='D38'; Value='1.06'}
    @{Cell='E38'; Value='  +6.44%  '}
    @{Cell='D39'; Value='90.77'}
    @{Cell='E39'; Value='  -1.22%  '}
    @{Cell='D40'; Value='1.344.97'}
    @{Cell='E40'; Value='  +2.38%  '}
    @{Cell='E41'; Value='  +0.39%  '}
    @{Cell='D42'; Value='15.03'}
    @{Cell='E42'; Value='  +3.68%  '}
    @{Cell='D43'; Value='2.30'}
    @{Cell='E43'; Value='  +0.57%  '}
    @{Cell='E44'; Value='  -2.71%  '}
    @{Cell='D45'; Value='2.75'}
    @{Cell='E45'; Value='  -0.28%  '}
    @{Cell='D46'; Value='6.29'}
    @{Cell='E46'; Value='  +0.07%  '}
    @{Cell='D47'; Value='0.0523'}
    @{Cell='D48'; Value='2.020.03'}
    @{Cell='E48'; Value='  +1.46%  '}
    @{Cell='B49'; Value='THORChain'}
    @{Cell='C49'; Value='https://coinranking.com/coin/ybmU-kKU+thorchain-rune'}
    @{Cell='D49'; Value='3.40'}
    @{Cell='E49'; Value='  +20.09%  '}
    @{Cell='B50'; Value='PaxDollar'}
    @{Cell='C50'; Value='https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'}
    @{Cell='D50'; Value='1.01'}
    @{Cell='E50'; Value='  -0.03%  '}
    @{Cell='D51'; Value='0.0665'}
    @{Cell='E51'; Value='  +1.77%  '}
)

foreach ($edit in $edits) {
    $rng = $ws.Range($edit.Cell)
    # Numeric-looking values (e.g. "231.79") would otherwise be auto-coerced to
    # a Number by COM's type inference; force Text so the stored cell type stays
    # a string, matching the workbook's existing inline-string convention. Restore
    # the original style afterwards so NumberFormat="@" leaves no visible trace.
    $savedStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $edit.Value
    $rng.Style = $savedStyle
}
